$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs) ---

# A8: "Volume 31   Number  46" -> "...47" (the "46" run becomes "47")
$rngA8 = $ws.Range("A8")
$volChars = $rngA8.Characters(21, 2)
$volChars.Text = "47"
$volChars.Font.Name = "Andale WT"
$volChars.Font.Size = 10

# C9: "Report Covering the Week  11/11/2024  Through  11/17/2024"
#     -> "...11/18/2024  Through  11/24/2024"
$rngC9 = $ws.Range("C9")
$dateChars1 = $rngC9.Characters(27, 10)
$dateChars1.Text = "11/18/2024"
$dateChars1.Font.Name = "Andale WT"
$dateChars1.Font.Size = 10
$dateChars2 = $rngC9.Characters(48, 10)
$dateChars2.Text = "11/24/2024"
$dateChars2.Font.Name = "Andale WT"
$dateChars2.Font.Size = 10

# --- Weekly crime-statistics table updates (rows 15-28) ---

# Cells that flip from the "n/a" text placeholder to a real number:
# use Value + NumberFormat (reuses the existing numeric style).
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"

# Cells that flip from a real number back to the "n/a" text placeholder:
# force Text format + string value, then copy the number-format/style from
# a known "n/a"-styled donor cell (C14) so the cell reuses the original style.
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

# Plain numeric value updates (style unchanged):
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 75
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -50
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 41.666666666666
$ws.Range("L16").Value = 6.25
$ws.Range("M16").Value = -16.393442622950
$ws.Range("N16").Value = -79.518072289156
$ws.Range("C17").Value = 3
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 143
$ws.Range("K17").Value = 25.438596491228
$ws.Range("L17").Value = 44.444444444444
$ws.Range("M17").Value = 83.333333333333
$ws.Range("N17").Value = -29.556650246305
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 51
$ws.Range("K18").Value = 30.769230769230
$ws.Range("L18").Value = -1.923076923076
$ws.Range("M18").Value = -3.773584905660
$ws.Range("N18").Value = -86.787564766839
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -17.647058823529
$ws.Range("I19").Value = 146
$ws.Range("J19").Value = 176
$ws.Range("K19").Value = -17.045454545454
$ws.Range("L19").Value = -5.806451612903
$ws.Range("M19").Value = 37.735849056603
$ws.Range("N19").Value = -19.337016574585
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 36.842105263157
$ws.Range("L20").Value = -7.142857142857
$ws.Range("M20").Value = 108
$ws.Range("N20").Value = -82.312925170068
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 34
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = 9.677419354838
$ws.Range("I21").Value = 450
$ws.Range("J21").Value = 407
$ws.Range("K21").Value = 10.565110565110
$ws.Range("L21").Value = 6.888361045130
$ws.Range("M21").Value = 37.614678899082
$ws.Range("N21").Value = -66.165413533834
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 15
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = -6.25
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 50
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 24
$ws.Range("H24").Value = 33.333333333333
$ws.Range("I24").Value = 458
$ws.Range("J24").Value = 429
$ws.Range("K24").Value = 6.759906759906
$ws.Range("L24").Value = -6.339468302658
$ws.Range("M24").Value = 65.942028985507
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 234
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = 5.882352941176
$ws.Range("L25").Value = 10.900473933649
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 239
$ws.Range("J26").Value = 202
$ws.Range("K26").Value = 18.316831683168
$ws.Range("L26").Value = 20.100502512562
$ws.Range("M26").Value = -23.397435897435
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 42.857142857142
$ws.Range("H28").Value = -100

$excel.CutCopyMode = 0
